# "Natmi following Dr Hou advice"
#
# The NATMI ligand-receptor table (Lama2 -> Rpsa) is recomputed after adding
# a new "ECs" (endothelial cells) sending/target cluster to the analysis.
# Previously the sheet only had the FAPs/sCs combinations (6 data rows); now
# every one of the 3 clusters (ECs, FAPs, sCs) is crossed with every other
# cluster as a target, giving the full 3x3 = 9 data rows, and all of the
# underlying statistics are recalculated with the new cluster included.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Row number, followed by the 20 column values (A:T) for that row.
$rowsData = @(
  @(2,  @("ECs",  "Lama2", "Rpsa", "ECs",  2, 0.6666666666666666, 1.471482,           4.414446,           0.004946458467382327, 0.004946458467382326, 3, 1, 106.5625623333333, 319.687687, 0.4373345410925676, 0.4373345410925676, 156.804892347378,   1411.244031126402,   0.002163257143866095,   0.002163257143866095)),
  @(3,  @("ECs",  "Lama2", "Rpsa", "FAPs", 2, 0.6666666666666666, 1.471482,           4.414446,           0.004946458467382327, 0.004946458467382326, 3, 1, 102.9000496666667, 308.700149, 0.4223035277493257, 0.4223035277493257, 151.415570883606,   1362.740137952454,   0.00208890686064108,    0.002088906860641079)),
  @(4,  @("ECs",  "Lama2", "Rpsa", "sCs",  2, 0.6666666666666666, 1.471482,           4.414446,           0.004946458467382327, 0.004946458467382326, 3, 1, 34.20111066666666, 102.603332, 0.1403619311581067, 0.1403619311581067, 50.326318726008,    452.9368685340719,   0.0006942944628751521,  0.000694294462875152)),
  @(5,  @("FAPs", "Lama2", "Rpsa", "ECs",  3, 1,                  264.7713316666666,  794.313995,         0.8900417371348598,   0.8900417371348596,   3, 1, 106.5625623333333, 319.687687, 0.4373345410925676, 0.4373345410925676, 28214.71153480884,  253932.4038132796,  0.3892459946631056,    0.3892459946631055)),
  @(6,  @("FAPs", "Lama2", "Rpsa", "FAPs", 3, 1,                  264.7713316666666,  794.313995,         0.8900417371348598,   0.8900417371348596,   3, 1, 102.9000496666667, 308.700149, 0.4223035277493257, 0.4223035277493257, 27244.98317880947,  245204.8486092853,  0.3758677654361893,    0.3758677654361892)),
  @(7,  @("FAPs", "Lama2", "Rpsa", "sCs",  3, 1,                  264.7713316666666,  794.313995,         0.8900417371348598,   0.8900417371348596,   3, 1, 34.20111066666666, 102.603332, 0.1403619311581067, 0.1403619311581067, 9055.47361569237,   81499.26254123134,  0.1249279770355649,    0.1249279770355649)),
  @(8,  @("sCs",  "Lama2", "Rpsa", "ECs",  3, 1,                  31.239114,          93.717342,          0.105011804397758,    0.105011804397758,    3, 1, 106.5625623333333, 319.687687, 0.4373345410925676, 0.4373345410925676, 3328.920032863106,  29960.28029576795,  0.04592528928559599,   0.04592528928559597)),
  @(9,  @("sCs",  "Lama2", "Rpsa", "FAPs", 3, 1,                  31.239114,          93.717342,          0.105011804397758,    0.105011804397758,    3, 1, 102.9000496666667, 308.700149, 0.4223035277493257, 0.4223035277493257, 3214.506382142662,  28930.55743928396,  0.04434685545249538,   0.04434685545249537)),
  @(10, @("sCs",  "Lama2", "Rpsa", "sCs",  3, 1,                  31.239114,          93.717342,          0.105011804397758,    0.105011804397758,    3, 1, 34.20111066666666, 102.603332, 0.1403619311581067, 0.1403619311581067, 1068.412395042616,  9615.711555383543,  0.01473965965966668,   0.01473965965966668))
)

foreach ($rowdef in $rowsData) {
  $rnum = $rowdef[0]
  $vals = $rowdef[1]
  for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + $rnum).Value = $vals[$i]
  }
}
